$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Peer review: reviewer name "Mina" added for every existing test case row
# (M2:M14), mirroring the "Reviewed by" column already on the sheet.
$ws.Range("M2:M14").Value = "Mina"

# Scroll the sheet so row 13 / column H is the top-left visible cell, and
# leave the M2:M14 range selected, as it was when the review was recorded.
$win = $excel.ActiveWindow
$win.ScrollRow = 13
$win.ScrollColumn = 8

$ws.Range("M2:M14").Select()
